# Refresh the coin Price (D) and Volume(1h) (E) columns with newly scraped
# values, row by row, matching the upstream GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.093.07'
$ws.Range("E2").Value = '  -1.31%  '

$ws.Range("D3").Value = '2.688.20'
$ws.Range("E3").Value = '  -1.85%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.88'
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.15'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.51%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("E8").Value = '  -0.70%  '

$ws.Range("E9").Value = '  -2.74%  '

$ws.Range("E10").Value = '  -2.05%  '

$ws.Range("E11").Value = '  -3.41%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.41'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.99%  '

$ws.Range("D13").Value = '3.162.37'
$ws.Range("E13").Value = '  -1.92%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.58'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.56%  '

$ws.Range("D15").Value = '62.988.32'
$ws.Range("E15").Value = '  -1.20%  '

$ws.Range("E16").Value = '  -1.42%  '

$ws.Range("D17").Value = '2.688.44'
$ws.Range("E17").Value = '  -1.97%  '

$ws.Range("E18").Value = '  -1.40%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.63'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '346.86'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.99%  '

$ws.Range("E21").Value = '  -4.40%  '

$ws.Range("E22").Value = '  +0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.513'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.70%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.43'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.42%  '

$ws.Range("E25").Value = '  -0.75%  '

$ws.Range("E26").Value = '  +0.28%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.26'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.95%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.44'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +8.68%  '

$ws.Range("E29").Value = '  -4.65%  '

$ws.Range("E30").Value = '  +0.69%  '

$ws.Range("E31").Value = '  -0.42%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '165.34'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.86%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.97'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.42%  '

$ws.Range("E34").Value = '  +1.30%  '

$ws.Range("E35").Value = '  +0.04%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '19.56'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.65%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.80'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.14%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '359.31'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.47%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.45'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.66%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.963'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.46%  '

$ws.Range("E41").Value = '  -1.74%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '38.53'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.17%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.17'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.08%  '

$ws.Range("E44").Value = '  -2.75%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0565'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.19%  '

$ws.Range("E46").Value = '  -0.38%  '

$ws.Range("E47").Value = '  -0.03%  '

$ws.Range("E48").Value = '  -0.14%  '

$ws.Range("E49").Value = '  -2.73%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0976'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.79%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '130.00'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.29%  '
